$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.450.62'
$ws.Range("E2").Value = '  +7.24%  '
$ws.Range("D3").Value = '1.817.55'
$ws.Range("E3").Value = '  +6.42%  '
$ws.Range("D4").Value = "'" + '1.003'
$ws.Range("E4").Value = '  +0.53%  '
$ws.Range("D5").Value = "'" + '343.73'
$ws.Range("E5").Value = '  +4.78%  '
$ws.Range("D6").Value = "'" + '0.9990'
$ws.Range("E6").Value = '  +0.49%  '
$ws.Range("D7").Value = "'" + '0.3836'
$ws.Range("E7").Value = '  +4.57%  '
$ws.Range("D8").Value = "'" + '50.29'
$ws.Range("E8").Value = '  +4.60%  '
$ws.Range("D9").Value = "'" + '0.3526'
$ws.Range("E9").Value = '  +7.26%  '
$ws.Range("D10").Value = "'" + '1.236'
$ws.Range("E10").Value = '  +5.90%  '
$ws.Range("D11").Value = "'" + '0.07836'
$ws.Range("E11").Value = '  +6.66%  '
$ws.Range("D12").Value = "'" + '0.9998'
$ws.Range("E12").Value = '  +0.58%  '
$ws.Range("D13").Value = "'" + '22.48'
$ws.Range("E13").Value = '  +12.43%  '
$ws.Range("D14").Value = "'" + '6.647'
$ws.Range("E14").Value = '  +7.31%  '
$ws.Range("D15").Value = "'" + '7.262'
$ws.Range("E15").Value = '  +6.66%  '
$ws.Range("D16").Value = '1.818.99'
$ws.Range("E16").Value = '  +7.32%  '
$ws.Range("D17").Value = "'" + '0.00001128'
$ws.Range("E17").Value = '  +5.11%  '
$ws.Range("D18").Value = "'" + '0.06770'
$ws.Range("E18").Value = '  +2.80%  '
$ws.Range("D19").Value = "'" + '87.52'
$ws.Range("E19").Value = '  +8.48%  '
$ws.Range("D20").Value = "'" + '0.9993'
$ws.Range("E20").Value = '  +0.57%  '
$ws.Range("D21").Value = "'" + '17.95'
$ws.Range("E21").Value = '  +11.15%  '
$ws.Range("D22").Value = "'" + '6.576'
$ws.Range("E22").Value = '  +8.51%  '
$ws.Range("D23").Value = "'" + '13.24'
$ws.Range("E23").Value = '  +1.09%  '
$ws.Range("D24").Value = '27.487.13'
$ws.Range("E24").Value = '  +7.40%  '
$ws.Range("D25").Value = "'" + '2.465'
$ws.Range("E25").Value = '  +1.16%  '
$ws.Range("D26").Value = "'" + '2.710'
$ws.Range("E26").Value = '  +9.24%  '
$ws.Range("D27").Value = "'" + '22.10'
$ws.Range("E27").Value = '  +15.38%  '
$ws.Range("D28").Value = "'" + '1.522'
$ws.Range("E28").Value = '  +19.24%  '
$ws.Range("D29").Value = "'" + '153.73'
$ws.Range("E29").Value = '  +2.45%  '
$ws.Range("D30").Value = '2.021.78'
$ws.Range("E30").Value = '  +7.34%  '
$ws.Range("D31").Value = "'" + '137.57'
$ws.Range("E31").Value = '  +7.00%  '
$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").Value = "'" + '6.422'
$ws.Range("E32").Value = '  +7.33%  '
$ws.Range("B33").Value = 'HuobiToken'
$ws.Range("C33").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D33").Value = "'" + '4.176'
$ws.Range("E33").Value = '  +2.18%  '
$ws.Range("B34").Value = 'Aptos'
$ws.Range("C34").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D34").Value = "'" + '13.84'
$ws.Range("E34").Value = '  +8.79%  '
$ws.Range("B35").Value = 'Stellar'
$ws.Range("C35").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D35").Value = "'" + '0.08829'
$ws.Range("E35").Value = '  +4.09%  '
$ws.Range("D36").Value = "'" + '1.728'
$ws.Range("E36").Value = '  +2.38%  '
$ws.Range("D37").Value = "'" + '5.679'
$ws.Range("E37").Value = '  +7.12%  '
$ws.Range("B38").Value = 'Hedera'
$ws.Range("C38").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D38").Value = "'" + '0.06568'
$ws.Range("E38").Value = '  +5.55%  '
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").Value = "'" + '0.02430'
$ws.Range("E39").Value = '  +7.07%  '
$ws.Range("B40").Value = 'Algorand'
$ws.Range("C40").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D40").Value = "'" + '0.2271'
$ws.Range("E40").Value = '  +7.01%  '
$ws.Range("B41").Value = 'TheSandbox'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D41").Value = "'" + '0.6911'
$ws.Range("E41").Value = '  +13.67%  '
$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").Value = "'" + '9.062'
$ws.Range("E42").Value = '  +6.92%  '
$ws.Range("E43").Value = '  -0.98%  '
$ws.Range("D44").Value = "'" + '15.11'
$ws.Range("E44").Value = '  +7.53%  '
$ws.Range("B45").Value = 'Decentraland'
$ws.Range("C45").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D45").Value = "'" + '0.6573'
$ws.Range("E45").Value = '  +12.59%  '
$ws.Range("B46").Value = 'Frax'
$ws.Range("C46").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D46").Value = "'" + '0.9987'
$ws.Range("E46").Value = '  +0.54%  '
$ws.Range("D47").Value = "'" + '3.972'
$ws.Range("E47").Value = '  +3.86%  '
$ws.Range("D48").Value = "'" + '2.183'
$ws.Range("E48").Value = '  +8.96%  '
$ws.Range("D49").Value = "'" + '133.23'
$ws.Range("E49").Value = '  +6.11%  '
$ws.Range("D50").Value = "'" + '0.07367'
$ws.Range("E50").Value = '  +2.08%  '
$ws.Range("D51").Value = "'" + '81.09'
$ws.Range("E51").Value = '  +6.72%  '
